$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=31; A=9861; B=43243; C=1; D=7; E=0; F=0.8; G=24; H=6; I="highpass"; J="lots of units, some burst may split units to 2!" },
    @{ Row=32; A=9861; B=43243; C=2; D=7; E=0; F=0.8; G=24; H=6; I="highpass"; J="no isolated cells" },
    @{ Row=33; A=9861; B=43243; C=3; D=7; E=0; F=0.8; G=24; H=6; I="highpass"; J=$null },
    @{ Row=34; A=9861; B=43243; C=4; D=7; E=0; F=0.8; G=24; H=6; I="highpass"; J=$null },
    @{ Row=35; A=9861; B=43244; C=1; D=7; E=0; F=0.8; G=24; H=6; I="highpass"; J=$null },
    @{ Row=36; A=9861; B=43244; C=1; D=6; E=0; F=0.8; G=24; H=6; I="highpass"; J=$null },
    @{ Row=37; A=9861; B=43244; C=2; D=7; E=0; F=0.8; G=24; H=6; I="highpass"; J=$null },
    @{ Row=38; A=9861; B=43244; C=3; D=7; E=0; F=0.8; G=24; H=6; I="highpass"; J=$null },
    @{ Row=39; A=9861; B=43244; C=4; D=7; E=0; F=0.8; G=24; H=6; I="highpass"; J="no cells" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item(30, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = $r.B

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    if ($r.J -ne $null) {
        $ws.Cells.Item($row, 10).Value = $r.J
    }
}

$excel.CutCopyMode = $false

$ws.Range("C40").Select()
